$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (I1, J1) - copy formatting from the existing header cell (H1)
# so they pick up the same bold/border/alignment style, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data cells on row 2 (plain numbers, no special style)
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 9
